$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 21.1360299999124
$ws.Range("E2").Value = [double]"-7.720532990428995E-07"
$ws.Range("F2").Value = 0.2417821421664148
$ws.Range("G2").Value = 3693.176271636694
$ws.Range("H2").Value = 0.5722995179579016
$ws.Range("B3").Value = 21.87472654997158
$ws.Range("C3").Value = 5
$ws.Range("E3").Value = [double]"3.86029E-07"
$ws.Range("F3").Value = 0.2821717707289705
$ws.Range("G3").Value = 3784.657714280826
$ws.Range("H3").Value = 0.577984277611966
$ws.Range("B4").Value = 22.61282812997593
$ws.Range("E4").Value = [double]"3.860242990428994E-07"
$ws.Range("F4").Value = 0.3023839000805217
$ws.Range("G4").Value = 3935.765975215543
$ws.Range("H4").Value = 0.5745470709481789
$ws.Range("B5").Value = 23.31729931995427
$ws.Range("F5").Value = 0.3154641296999544
$ws.Range("G5").Value = 4080.704025106822
$ws.Range("H5").Value = 0.5714038356247579
$ws.Range("B6").Value = 24.10071432995696
$ws.Range("F6").Value = 0.3200592976877065
$ws.Range("G6").Value = 4299.186081075621
$ws.Range("H6").Value = 0.5605878386154237
$ws.Range("B7").Value = 24.98359055991395
$ws.Range("D7").Value = 5
$ws.Range("E7").Value = 27.21668751492523
$ws.Range("F7").Value = 0.3667547375483475
$ws.Range("G7").Value = 4487.59752519654
$ws.Range("H7").Value = 0.5567252949855338
$ws.Range("B8").Value = 25.97364393996251
$ws.Range("D8").Value = 7
$ws.Range("E8").Value = 31.96831289309676
$ws.Range("F8").Value = 0.3670829004161846
$ws.Range("G8").Value = 4629.628121691282
$ws.Range("H8").Value = 0.5610308918391894
$ws.Range("B9").Value = 27.19861199996599
$ws.Range("D9").Value = 13
$ws.Range("E9").Value = 43.58714678157924
$ws.Range("F9").Value = 0.386164443454968
$ws.Range("G9").Value = 4729.689828867503
$ws.Range("H9").Value = 0.5750612193205604
$ws.Range("C10").Value = 13
$ws.Range("D10").Value = 24
$ws.Range("E10").Value = 170.0640501765404
$ws.Range("F10").Value = 0.3891590539204214
$ws.Range("G10").Value = 4869.217031722656
$ws.Range("H10").Value = 0.5710469631732992
$ws.Range("B11").Value = 28.17733188996433
$ws.Range("C11").Value = 7
$ws.Range("D11").Value = 28
$ws.Range("E11").Value = 186.6001118196207
$ws.Range("F11").Value = 0.3810818663094769
$ws.Range("G11").Value = 5017.485049955806
$ws.Range("H11").Value = 0.5615827772164964
$ws.Range("B12").Value = 28.57310893998299
$ws.Range("C12").Value = 25
$ws.Range("D12").Value = 41
$ws.Range("E12").Value = 323.2522956176013
$ws.Range("F12").Value = 0.3371212568904194
$ws.Range("G12").Value = 5231.346959647781
$ws.Range("H12").Value = 0.5461902863714237
$ws.Range("B13").Value = 28.96257434995407
$ws.Range("C13").Value = 6
$ws.Range("D13").Value = 39
$ws.Range("E13").Value = 326.2799606409457
$ws.Range("F13").Value = 0.3247564344106765
$ws.Range("G13").Value = 5371.465998473741
$ws.Range("H13").Value = 0.5391931059078389
$ws.Range("B14").Value = 29.29447642992334
$ws.Range("C14").Value = 7
$ws.Range("D14").Value = 31
$ws.Range("E14").Value = 230.9913779262596
$ws.Range("F14").Value = 0.3245582902341944
$ws.Range("G14").Value = 5508.782784766356
$ws.Range("H14").Value = 0.5317776644040578
$ws.Range("B15").Value = 29.23475238995202
$ws.Range("C15").Value = 4
$ws.Range("E15").Value = 146.5967438839343
$ws.Range("F15").Value = 0.3136375020813455
$ws.Range("G15").Value = 5556.338554246508
$ws.Range("H15").Value = 0.5261513873665773
$ws.Range("B16").Value = 29.17600789995107
$ws.Range("C16").Value = 12
$ws.Range("E16").Value = 173.5731612614937
$ws.Range("F16").Value = 0.2919396131570683
$ws.Range("G16").Value = 5630.543830587956
$ws.Range("H16").Value = 0.5181738883099047
$ws.Range("B17").Value = 29.0963684099763
$ws.Range("C17").Value = 7
$ws.Range("D17").Value = 24
$ws.Range("E17").Value = 168.3629475157368
$ws.Range("F17").Value = 0.2816061528257992
$ws.Range("G17").Value = 5663.074400323652
$ws.Range("H17").Value = 0.5137910320993382
$ws.Range("B18").Value = 28.96841058994999
$ws.Range("C18").Value = 3
$ws.Range("E18").Value = 114.5281331950805
$ws.Range("F18").Value = 0.2819566803806406
$ws.Range("G18").Value = 5737.793202661475
$ws.Range("H18").Value = 0.5048702448271052
$ws.Range("B19").Value = 28.81461644995066
$ws.Range("D19").Value = 5
$ws.Range("E19").Value = 5.621839967628925
$ws.Range("F19").Value = 0.2819936382376502
$ws.Range("G19").Value = 5739.149036941112
$ws.Range("H19").Value = 0.5020712350294436
$ws.Range("B20").Value = 28.30960012995091
$ws.Range("F20").Value = 0.2781906552442441
$ws.Range("G20").Value = 5683.490657990042
$ws.Range("H20").Value = 0.4981023429703774
$ws.Range("B21").Value = 27.7800253099503
$ws.Range("F21").Value = 0.257267396543277
$ws.Range("G21").Value = 5620.248905096278
$ws.Range("H21").Value = 0.4942846087254283
$ws.Range("B22").Value = 27.19968143994861
$ws.Range("F22").Value = 0.2457077183393178
$ws.Range("G22").Value = 5510.804137556207
$ws.Range("H22").Value = 0.4935700990456619
$ws.Range("B23").Value = 26.57332149995102
$ws.Range("E23").Value = 0.1367933993437299
$ws.Range("F23").Value = 0.245328087940182
$ws.Range("G23").Value = 5363.870386796181
$ws.Range("H23").Value = 0.4954131920369382
$ws.Range("B24").Value = 25.90453326995338
$ws.Range("E24").Value = 0.153541385
$ws.Range("F24").Value = 0.2453267310236783
$ws.Range("G24").Value = 5181.798617837004
$ws.Range("H24").Value = 0.4999139329881272
$ws.Range("B25").Value = 25.1162635900469
$ws.Range("F25").Value = 0.2417506312467248
$ws.Range("G25").Value = 5036.97773473263
$ws.Range("H25").Value = 0.4986375742115546
$ws.Range("B26").Value = 24.31404712995696
$ws.Range("F26").Value = 0.2386332399090469
$ws.Range("G26").Value = 4874.709737789721
$ws.Range("H26").Value = 0.4987793825234276
$ws.Range("B27").Value = 23.51812297995581
$ws.Range("E27").Value = 0.153541385
$ws.Range("F27").Value = 0.2333722968876059
$ws.Range("G27").Value = 4746.364448290897
$ws.Range("H27").Value = 0.4954976221521373
$ws.Range("B28").Value = 22.72592122995475
$ws.Range("F28").Value = 0.2237523068202285
$ws.Range("G28").Value = 4537.298804107789
$ws.Range("H28").Value = 0.5008689577459636
$ws.Range("B29").Value = 21.93048930995987
$ws.Range("C29").Value = 3
$ws.Range("D29").Value = 3
$ws.Range("E29").Value = 8.289760027
$ws.Range("F29").Value = 0.2229922774237511
$ws.Range("G29").Value = 4374.471505543002
$ws.Range("H29").Value = 0.5013288869791747
